$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: mark item as "Desenvolvido" (was "Em análise") and restyle like the
#     other already-resolved rows (green fill, style used by row 19/20). ---
$ws.Range("A19:D19").Copy() | Out-Null
$ws.Range("A17:D17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value = "Desenvolvido"

# --- Row 18: mark item as "Analisado" / "Não procede" and restyle to match. ---
$ws.Range("A20:D20").Copy() | Out-Null
$ws.Range("A18:D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = "Não procede"
$ws.Range("C18").Value = "Analisado"

# --- Row 26: same treatment as row 18. ---
$ws.Range("A20:D20").Copy() | Out-Null
$ws.Range("A26:D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "Não procede"
$ws.Range("C26").Value = "Analisado"

$excel.CutCopyMode = 0

# --- Update the active selection on the sheet. ---
$ws.Range("A22").Select()
